$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "G2"
$ws.Range("B3").Value = "Workout"
$ws.Range("C3").Value = "Daily"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 45907
$ws.Range("E3").Style = $ws.Range("E2").Style
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("F3").Value = 36
